$wb = $excel.ActiveWorkbook

# ============================================================
# Step 1: reorganise sheets
#   before: ODI Batting, ODI Bowling
#   after:  Player Info, ODI Batting, ODI Bowling, ODI Batting Extra
# ============================================================
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

# re-fetch by name: inserting a sheet shifts indices, and object refs
# captured before the insert can resolve against the stale position
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$battingExtra = $wb.Worksheets.Add($null, $bowlingSheet)
$battingExtra.Name = "ODI Batting Extra"

# ============================================================
# Step 2: "ODI Batting" -- MATCH_CARD_LINK column becomes MATCH_CODE
# ============================================================
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value2 = "MATCH_CODE"
$battingSheet.Range("D2:D61").NumberFormat = "@"
$battingSheet.Range("D2").Value2 = '3128'
$battingSheet.Range("D3").Value2 = '3174'
$battingSheet.Range("D4").Value2 = '3175'
$battingSheet.Range("D5").Value2 = '3176'
$battingSheet.Range("D6").Value2 = '3177'
$battingSheet.Range("D7").Value2 = '3178'
$battingSheet.Range("D8").Value2 = '3196'
$battingSheet.Range("D9").Value2 = '3197'
$battingSheet.Range("D10").Value2 = '3198'
$battingSheet.Range("D11").Value2 = '3201'
$battingSheet.Range("D12").Value2 = '3220'
$battingSheet.Range("D13").Value2 = '3232'
$battingSheet.Range("D14").Value2 = '3267'
$battingSheet.Range("D15").Value2 = '3274'
$battingSheet.Range("D16").Value2 = '3277'
$battingSheet.Range("D17").Value2 = '3282'
$battingSheet.Range("D18").Value2 = '3287'
$battingSheet.Range("D19").Value2 = '3288'
$battingSheet.Range("D20").Value2 = '3289'
$battingSheet.Range("D21").Value2 = '3290'
$battingSheet.Range("D22").Value2 = '3292'
$battingSheet.Range("D23").Value2 = '3330'
$battingSheet.Range("D24").Value2 = '3352'
$battingSheet.Range("D25").Value2 = '3358'
$battingSheet.Range("D26").Value2 = '3372'
$battingSheet.Range("D27").Value2 = '3374'
$battingSheet.Range("D28").Value2 = '3379'
$battingSheet.Range("D29").Value2 = '3383'
$battingSheet.Range("D30").Value2 = '3394'
$battingSheet.Range("D31").Value2 = '3410'
$battingSheet.Range("D32").Value2 = '3411'
$battingSheet.Range("D33").Value2 = '3413'
$battingSheet.Range("D34").Value2 = '3439'
$battingSheet.Range("D35").Value2 = '3441'
$battingSheet.Range("D36").Value2 = '3443'
$battingSheet.Range("D37").Value2 = '3484'
$battingSheet.Range("D38").Value2 = '3497'
$battingSheet.Range("D39").Value2 = '3499'
$battingSheet.Range("D40").Value2 = '3500'
$battingSheet.Range("D41").Value2 = '3506'
$battingSheet.Range("D42").Value2 = '3514'
$battingSheet.Range("D43").Value2 = '3531'
$battingSheet.Range("D44").Value2 = '3532'
$battingSheet.Range("D45").Value2 = '3574'
$battingSheet.Range("D46").Value2 = '3575'
$battingSheet.Range("D47").Value2 = '3585'
$battingSheet.Range("D48").Value2 = '3677'
$battingSheet.Range("D49").Value2 = '3679'
$battingSheet.Range("D50").Value2 = '3681'
$battingSheet.Range("D51").Value2 = '3713'
$battingSheet.Range("D52").Value2 = '3715'
$battingSheet.Range("D53").Value2 = '3717'
$battingSheet.Range("D54").Value2 = '3720'
$battingSheet.Range("D55").Value2 = '3801'
$battingSheet.Range("D56").Value2 = '3802'
$battingSheet.Range("D57").Value2 = '3803'
$battingSheet.Range("D58").Value2 = '3837'
$battingSheet.Range("D59").Value2 = '3838'
$battingSheet.Range("D60").Value2 = '3973'
$battingSheet.Range("D61").Value2 = '3975'

# two rows ("did not bat") have a stray empty INNING_NUMBER cell that
# should not be present at all
$battingSheet.Range("B16").ClearContents()
$battingSheet.Range("B55").ClearContents()

# ============================================================
# Step 3: "ODI Bowling" -- MATCH_CARD_LINK column becomes MATCH_CODE
# ============================================================
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value2 = "MATCH_CODE"
$bowlingSheet.Range("B2:B3").NumberFormat = "@"
$bowlingSheet.Range("B2").Value2 = '3713'
$bowlingSheet.Range("B3").Value2 = '3715'

# ============================================================
# Step 4: populate the new "Player Info" sheet
# ============================================================
$playerInfo = $wb.Worksheets.Item("Player Info")
$playerInfo.Range("A1:D2").NumberFormat = "@"
$playerInfo.Range("A1").Value2 = 'ID'
$playerInfo.Range("B1").Value2 = 'NAME'
$playerInfo.Range("C1").Value2 = 'BATTING_HAND'
$playerInfo.Range("D1").Value2 = 'BOWL_STYLE'
$playerInfo.Range("A2").Value2 = '3798'
$playerInfo.Range("B2").Value2 = 'Asad Shafiq'
$playerInfo.Range("C2").Value2 = 'Right Handed'
$playerInfo.Range("D2").Value2 = 'Right Arm Leg Break'
$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# ============================================================
# Step 5: populate the new "ODI Batting Extra" sheet
# ============================================================
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$battingExtra.Range("A2:A21").NumberFormat = "@"
$battingExtra.Range("C2:C21").NumberFormat = "@"
$battingExtra.Range("D2:D21").NumberFormat = "@"
$battingExtra.Range("E2:E21").NumberFormat = "@"
$battingExtra.Range("F2:F21").NumberFormat = "@"
$battingExtra.Range("A1").Value2 = 'MATCH_CODE'
$battingExtra.Range("B1").Value2 = 'BATTING_POSITION'
$battingExtra.Range("C1").Value2 = 'NUM_4'
$battingExtra.Range("D1").Value2 = 'NUM_6'
$battingExtra.Range("E1").Value2 = 'PERCENT_RUNS_OF_TOTAL'
$battingExtra.Range("F1").Value2 = 'MAN_OF_MATCH'
$battingExtra.Range("A2").Value2 = '3514'
$battingExtra.Range("B2").Value2 = 4
$battingExtra.Range("C2").Value2 = '3'
$battingExtra.Range("D2").Value2 = '0'
$battingExtra.Range("E2").Value2 = '24.85%'
$battingExtra.Range("F2").Value2 = 'NO'
$battingExtra.Range("A3").Value2 = '3531'
$battingExtra.Range("B3").Value2 = 5
$battingExtra.Range("C3").Value2 = '0'
$battingExtra.Range("D3").Value2 = '0'
$battingExtra.Range("F3").Value2 = 'NO'
$battingExtra.Range("A4").Value2 = '3532'
$battingExtra.Range("F4").Value2 = 'NO'
$battingExtra.Range("A5").Value2 = '3574'
$battingExtra.Range("F5").Value2 = 'NO'
$battingExtra.Range("A6").Value2 = '3575'
$battingExtra.Range("B6").Value2 = 4
$battingExtra.Range("C6").Value2 = '0'
$battingExtra.Range("D6").Value2 = '0'
$battingExtra.Range("E6").Value2 = '0.42%'
$battingExtra.Range("F6").Value2 = 'NO'
$battingExtra.Range("A7").Value2 = '3585'
$battingExtra.Range("B7").Value2 = 3
$battingExtra.Range("C7").Value2 = '0'
$battingExtra.Range("D7").Value2 = '0'
$battingExtra.Range("E7").Value2 = '0.56%'
$battingExtra.Range("F7").Value2 = 'NO'
$battingExtra.Range("A8").Value2 = '3677'
$battingExtra.Range("F8").Value2 = 'NO'
$battingExtra.Range("A9").Value2 = '3679'
$battingExtra.Range("B9").Value2 = 3
$battingExtra.Range("C9").Value2 = '2'
$battingExtra.Range("D9").Value2 = '1'
$battingExtra.Range("E9").Value2 = '13.49%'
$battingExtra.Range("F9").Value2 = 'NO'
$battingExtra.Range("A10").Value2 = '3681'
$battingExtra.Range("B10").Value2 = 3
$battingExtra.Range("C10").Value2 = '5'
$battingExtra.Range("D10").Value2 = '0'
$battingExtra.Range("E10").Value2 = '21.74%'
$battingExtra.Range("F10").Value2 = 'NO'
$battingExtra.Range("A11").Value2 = '3713'
$battingExtra.Range("B11").Value2 = 3
$battingExtra.Range("C11").Value2 = '0'
$battingExtra.Range("D11").Value2 = '0'
$battingExtra.Range("E11").Value2 = '2.00%'
$battingExtra.Range("F11").Value2 = 'NO'
$battingExtra.Range("A12").Value2 = '3715'
$battingExtra.Range("B12").Value2 = 4
$battingExtra.Range("C12").Value2 = '0'
$battingExtra.Range("D12").Value2 = '0'
$battingExtra.Range("E12").Value2 = '0.40%'
$battingExtra.Range("F12").Value2 = 'NO'
$battingExtra.Range("A13").Value2 = '3717'
$battingExtra.Range("B13").Value2 = 4
$battingExtra.Range("C13").Value2 = '1'
$battingExtra.Range("D13").Value2 = '0'
$battingExtra.Range("E13").Value2 = '6.32%'
$battingExtra.Range("F13").Value2 = 'NO'
$battingExtra.Range("A14").Value2 = '3720'
$battingExtra.Range("F14").Value2 = 'NO'
$battingExtra.Range("A15").Value2 = '3801'
$battingExtra.Range("B15").Value2 = 7
$battingExtra.Range("F15").Value2 = 'NO'
$battingExtra.Range("A16").Value2 = '3802'
$battingExtra.Range("F16").Value2 = 'NO'
$battingExtra.Range("A17").Value2 = '3803'
$battingExtra.Range("B17").Value2 = 3
$battingExtra.Range("C17").Value2 = '1'
$battingExtra.Range("D17").Value2 = '0'
$battingExtra.Range("E17").Value2 = '5.41%'
$battingExtra.Range("F17").Value2 = 'NO'
$battingExtra.Range("A18").Value2 = '3837'
$battingExtra.Range("F18").Value2 = 'NO'
$battingExtra.Range("A19").Value2 = '3838'
$battingExtra.Range("B19").Value2 = 4
$battingExtra.Range("C19").Value2 = '2'
$battingExtra.Range("D19").Value2 = '0'
$battingExtra.Range("E19").Value2 = '23.46%'
$battingExtra.Range("F19").Value2 = 'NO'
$battingExtra.Range("A20").Value2 = '3973'
$battingExtra.Range("F20").Value2 = 'NO'
$battingExtra.Range("A21").Value2 = '3975'
$battingExtra.Range("B21").Value2 = 4
$battingExtra.Range("C21").Value2 = '0'
$battingExtra.Range("D21").Value2 = '0'
$battingExtra.Range("E21").Value2 = '1.90%'
$battingExtra.Range("F21").Value2 = 'NO'
$headerRange2 = $battingExtra.Range("A1:F1")
$headerRange2.Font.Bold = $true
$headerRange2.Borders.LineStyle = 1
$headerRange2.HorizontalAlignment = -4108
$headerRange2.VerticalAlignment = -4160

Write-Host "edit complete"
